# Weekly refresh of the daily "Espinaca" price log: the existing data rows
# (2..24) are re-shuffled (each row's Fecha/Calidad/Volumen/Precio* /Origen
# values are reassigned to another row's prior values) while the
# market/category columns (A,B,C,E,F,G,H,N,Q,R) - which are identical for
# every row in this sheet - are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (the source row's D/I/J/K/L/M/O/P values get
# copied onto the destination row)
$rowMap = @{
    2  = 6
    3  = 9
    4  = 2
    5  = 16
    6  = 17
    7  = 22
    8  = 21
    9  = 18
    10 = 8
    11 = 12
    12 = 19
    13 = 14
    14 = 3
    15 = 13
    16 = 20
    17 = 4
    18 = 15
    19 = 24
    20 = 5
    21 = 11
    22 = 7
    23 = 10
    24 = 23
}

$cols = @('D', 'I', 'J', 'K', 'L', 'M', 'O', 'P')

# Snapshot every source row's values first (the remapping contains cycles,
# so cells must not be overwritten before they have been read).
$snapshot = @{}
foreach ($r in 2..24) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in 2..24) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
